$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated rand_digit values (column J) for rows 2-81
$randDigitValues = @(
    1, 8, 6, 5, 7, 7, 1, 2, 7, 1, 7, 1, 8, 5, 6, 4, 5, 8, 4, 7,
    5, 4, 3, 5, 6, 4, 6, 3, 6, 4, 3, 1, 6, 1, 8, 5, 7, 3, 3, 8,
    2, 5, 3, 5, 1, 4, 5, 4, 7, 6, 1, 1, 5, 5, 1, 5, 4, 4, 8, 1,
    8, 5, 7, 4, 6, 5, 2, 4, 7, 1, 6, 1, 6, 5, 8, 8, 4, 5, 5, 6
)

for ($i = 0; $i -lt $randDigitValues.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 10).Value = $randDigitValues[$i]
}
